# Replace the "Tony Tester" sample row (row 2) with the new "Rajesh" contact info.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Rajesh"
$ws.Range("B2").Value = "Rajesh@gmail.com"
$ws.Range("C2").Value = "Rajesh"

# Update the first data row (row 5) of the second table to the new "Rajesh 3" entry.
$ws.Range("A5").Value = "Rajesh 3"
$ws.Range("C5").Value = "Male 3"
$ws.Range("D5").Value = "A 3"

# Append three more data rows (6-8) to the second table, supporting 2-level reading.
# Fill column by column so the shared-string table is built in the same order
# Excel would use when the columns are typed/filled one after another.
$ws.Range("A6").Value = "Rajesh 4"
$ws.Range("A7").Value = "Rajesh 5"
$ws.Range("A8").Value = "Rajesh 6"

$ws.Range("B6").Value = 18
$ws.Range("B7").Value = 18
$ws.Range("B8").Value = 18

$ws.Range("C6").Value = "Male 4"
$ws.Range("C7").Value = "Male 5"
$ws.Range("C8").Value = "Male 6"

$ws.Range("D6").Value = "A 4"
$ws.Range("D7").Value = "A 5"
$ws.Range("D8").Value = "A 6"

# Match the author's final selection: D5:D8 with the active cell at D5.
$ws.Range("D5:D8").Select()
